$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week 23 (Z) and Week 24 (AA) columns are being added to the weekly report.
# Header row 1 holds the week numbers as text labels, like the existing D1:Y1 cells.
$ws.Range("Z1").Value = "'23"
$ws.Range("AA1").Value = "'24"
$ws.Range("Z1:AA1").Font.Bold = $true
$ws.Range("Z1:AA1").HorizontalAlignment = -4108

# Weekly counts per facility for weeks 23 and 24 (plus a couple of week 22 corrections).
$ws.Range("Z2").Value = 72
$ws.Range("AA2").Value = 67
$ws.Range("Z3").Value = 46
$ws.Range("Z4").Value = 1
$ws.Range("AA4").Value = 0
$ws.Range("Z5").Value = 4
$ws.Range("AA5").Value = 2
$ws.Range("Z6").Value = 95
$ws.Range("AA6").Value = 89
$ws.Range("Z7").Value = 34
$ws.Range("AA7").Value = 27
$ws.Range("Z8").Value = 42
$ws.Range("AA8").Value = 33
$ws.Range("Z9").Value = 1
$ws.Range("AA9").Value = 4
$ws.Range("Z10").Value = 2
$ws.Range("Z11").Value = 2
$ws.Range("Z12").Value = 3
$ws.Range("AA12").Value = 5
$ws.Range("Z13").Value = 2
$ws.Range("AA13").Value = 3
$ws.Range("Z14").Value = 1
$ws.Range("AA14").Value = 1
$ws.Range("Z15").Value = 2
$ws.Range("AA15").Value = 2
$ws.Range("Z16").Value = 2
$ws.Range("AA16").Value = 1
$ws.Range("Z17").Value = 5
$ws.Range("AA17").Value = 2
$ws.Range("Z18").Value = 1
$ws.Range("Z19").Value = 2
$ws.Range("Z21").Value = 3
$ws.Range("AA21").Value = 4
$ws.Range("Z22").Value = 7
$ws.Range("AA22").Value = 1
$ws.Range("AA23").Value = 1
$ws.Range("Z24").Value = 45
$ws.Range("AA24").Value = 47
$ws.Range("Z25").Value = 3
$ws.Range("AA25").Value = 1
$ws.Range("Z27").Value = 246
$ws.Range("AA27").Value = 229
$ws.Range("Y28").Value = 0
$ws.Range("Z28").Value = 0
$ws.Range("AA28").Value = 0
$ws.Range("Y29").Value = 32
$ws.Range("Z29").Value = 21
$ws.Range("AA29").Value = 20
$ws.Range("Z30").Value = 1
$ws.Range("AA30").Value = 2
$ws.Range("Z31").Value = 9
$ws.Range("AA31").Value = 7
$ws.Range("Z33").Value = 0
$ws.Range("AA33").Value = 0
$ws.Range("Z34").Value = 52
$ws.Range("AA34").Value = 38
$ws.Range("Y35").Value = 5
$ws.Range("Z35").Value = 6
$ws.Range("AA35").Value = 7
$ws.Range("Z36").Value = 10
$ws.Range("AA36").Value = 9
$ws.Range("Z37").Value = 73
$ws.Range("AA37").Value = 87
$ws.Range("Z39").Value = 14
$ws.Range("AA39").Value = 10
$ws.Range("Z40").Value = 52
$ws.Range("Y41").Value = 75
$ws.Range("Z41").Value = 33
$ws.Range("Y42").Value = 415
$ws.Range("Z42").Value = 219
$ws.Range("AA42").Value = 212
$ws.Range("Z43").Value = 70
$ws.Range("AA43").Value = 95
$ws.Range("Z44").Value = 146
$ws.Range("AA44").Value = 120
$ws.Range("Z45").Value = 1
$ws.Range("AA45").Value = 1
$ws.Range("Z46").Value = 98
$ws.Range("AA46").Value = 118
$ws.Range("Z47").Value = 6
$ws.Range("AA47").Value = 2
$ws.Range("Z48").Value = 0
$ws.Range("AA48").Value = 0
$ws.Range("AA49").Value = 12
$ws.Range("Y50").Value = 4
$ws.Range("Z50").Value = 1
$ws.Range("Z51").Value = 52
$ws.Range("AA51").Value = 39
$ws.Range("Z52").Value = 0
$ws.Range("AA52").Value = 0
$ws.Range("Z53").Value = 0
$ws.Range("AA53").Value = 0
$ws.Range("Z54").Value = 4
$ws.Range("AA54").Value = 6
$ws.Range("Z55").Value = 14
$ws.Range("AA55").Value = 0
$ws.Range("Z56").Value = 33
$ws.Range("AA56").Value = 37
